$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (shifts current rows 4-10 down to 5-11)
$ws.Rows.Item(4).Insert()

# Fill the new row 4: B4="dependency", C4="Module"
$ws.Range("B4").Value = "dependency"
$ws.Range("C4").Value = "Module"

# Insert a new row after the last data row (current row 11, which was row 10 before insert) to host Step3
$ws.Rows.Item(12).Insert()
$ws.Range("B12").Value = "Step3"
$ws.Range("C12").Value = "= doSomething ($Step2)"
